## Lab3 fix: remove the stray "[way, route] = dijkstra(paths, src, dest);"
## line (and the blank line above it) from the second for-loop, where
## `way`/`route` were already computed a couple of lines earlier from the
## Floyd results (`d`/`routes`) -- the extra Dijkstra call there was a bug
## left over from copy/pasting the first loop. Deleting the text the way a
## human editor would (select from just after "...in routes else [];"
## through the end of the erroneous statement's paragraph, then hit
## Delete) leaves Word's usual "_GoBack" bookmark behind at the edit point.

$d = $word.ActiveDocument

# Locate the unique anchor text that immediately precedes the blank line
# and the erroneous Dijkstra call (only occurs once in the document).
$anchor = $d.Content
$found = $anchor.Find.Execute("in routes else []", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "anchor text not found"
}

# Move to the end of the anchor's paragraph (end of the "route = ...;" line).
$anchor.Collapse(0)
[void]$anchor.Expand(4)            # wdParagraph
$afterAnchorPara = $anchor.End

# The blank paragraph right after it.
$blankPara = $d.Range($afterAnchorPara, $afterAnchorPara)
[void]$blankPara.Expand(4)         # wdParagraph

# The erroneous "[way, route] = dijkstra(paths, src, dest);" paragraph.
$badPara = $d.Range($blankPara.End, $blankPara.End)
[void]$badPara.Expand(4)           # wdParagraph

# Delete both paragraphs (including their paragraph marks) in one go so the
# following "print(...)" line moves up to take their place.
$deleteRange = $d.Range($blankPara.Start, $badPara.End)
$deleteRange.Delete()

# Word drops a "_GoBack" bookmark at the last edited spot.
$goBackRange = $d.Range($blankPara.Start, $blankPara.Start)
$d.Bookmarks.Add("_GoBack", $goBackRange)
